$d = $word.ActiveDocument

function Get-ParagraphAt($pos) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

# Locate the three consecutive list-item paragraphs being collapsed into one:
#   "Anvil alter cloth has stencil property of cull none but my nifdisplay is
#    definitely cull back"
#   "Make a new sub of indexedgeom and make it call the interleaved interface,
#    give it everything properly including short indexes"
#   "Do the animation back to transforms with a double buffer idea"
$findFirst = $d.Content
$findFirst.Find.Execute("Anvil alter cloth has s")
$firstPara = Get-ParagraphAt($findFirst.Start)

$findSecond = $d.Content
$findSecond.Find.Execute("Make a new sub of ")
$secondPara = Get-ParagraphAt($findSecond.Start)

$findThird = $d.Content
$findThird.Find.Execute("Do the animation back to transforms with a double buffer idea")
$thirdPara = Get-ParagraphAt($findThird.Start)

# Remove the 2nd and 3rd paragraphs entirely (delete from the end backwards so
# the earlier paragraph's position stays valid).
$thirdPara.Range.Delete()
$secondPara.Range.Delete()

# Rewrite the remaining (first) paragraph's text.
$rng = $firstPara.Range
$rng.End = $rng.End - 1
$rng.Text = "Do the animation double buffer idea"

# Append the trailing "?" as its own run (matching the target XML, which keeps
# it as a separate <w:r> even though formatting is identical to the run
# before it). Toggling a character property and back forces the engine to
# keep the inserted text as a distinct run instead of merging it into the
# adjacent, identically-formatted run.
$afterRng = $firstPara.Range
$afterRng.End = $afterRng.End - 1
$insertPos = $afterRng.End
$qRng = $d.Range($insertPos, $insertPos)
$qRng.InsertAfter("?")
$qMark = $d.Range($insertPos, $insertPos + 1)
$qMark.Font.Bold = $true
$qMark.Font.Bold = $false
